$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) cells - use a leading apostrophe to force text entry
# (these look like numbers/dates to Excel auto-detection), then reset the style
# so no extra number-format / quote-prefix styling is left on the cell.
$ws.Range("D2").Value = "'62.945.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.448.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'570.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'145.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'2.447.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Value = "'5.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.355"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'26.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Value = "'62.952.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'2.448.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'11.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'328.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'65.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'613.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'8.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Value = "'2.580.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'1.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Value = "'8.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'1.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.140"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Value = "'0.378"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'18.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'5.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'146.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'1.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Value = "'41.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'148.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'3.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'21.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Value = "'0.601"
$ws.Range("D50").Style = "Normal"

# Update Volume(1h) column (E) cells - plain text assignment is safe because
# the values carry padding spaces and a percent sign that keep Excel from
# reinterpreting them as numbers.
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  +5.04%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  +13.20%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("E43").Value = "  +4.58%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +0.56%  "
